$wb = $excel.ActiveWorkbook

# Update OFF sheet (row 3 cumulative stats, Week 15 logged / Week 16 simulated)
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B3").Value = 338
$wsOff.Range("C3").Value = 226
$wsOff.Range("D3").Value = 67
$wsOff.Range("E3").Value = 24
$wsOff.Range("F3").Value = 9

# Update DEF sheet (row 3 cumulative stats)
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B3").Value = 443
$wsDef.Range("C3").Value = 316
$wsDef.Range("D3").Value = 102
$wsDef.Range("E3").Value = 53
